$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 510.58334
$ws.Range("I4").Value = 312.7
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 312.7
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -198.7
$ws.Range("N4").Value = -1728
$ws.Range("H113").Value = 3404.4211
$ws.Range("I113").Value = 4398
$ws.Range("K113").Value = 4398
$ws.Range("M113").Value = -1144
$ws.Range("H129").Value = 1157.9
$ws.Range("I129").Value = 1548.8889
$ws.Range("J129").Value = 990.3333
$ws.Range("K129").Value = 4646.6667
$ws.Range("L129").Value = 2970.9999
$ws.Range("M129").Value = 353.3333000000002
$ws.Range("N129").Value = -12970.9999
$ws.Range("H134").Value = 38646.355
$ws.Range("J134").Value = 38646.355
$ws.Range("L134").Value = 38646.355
$ws.Range("N134").Value = -48786.355
$ws.Range("H136").Value = 38613.57
$ws.Range("J136").Value = 38613.57
$ws.Range("L136").Value = 38613.57
$ws.Range("N136").Value = -48813.57
$ws.Range("H138").Value = 2010.6912
$ws.Range("I138").Value = 2084.348
$ws.Range("J138").Value = 1973.0444
$ws.Range("K138").Value = 6253.044
$ws.Range("L138").Value = 5919.1332
$ws.Range("M138").Value = -1113.044
$ws.Range("N138").Value = -16199.1332

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 942.3611
$ws.Range("I2").Value = 887.0357
$ws.Range("J2").Value = 1136
$ws.Range("K2").Value = 887.0357
$ws.Range("L2").Value = 1136
$ws.Range("M2").Value = -774.0357
$ws.Range("N2").Value = -1362
$ws.Range("H26").Value = 1174.75
$ws.Range("I26").Value = 1174.75
$ws.Range("K26").Value = 1174.75
$ws.Range("M26").Value = -844.75
$ws.Range("H97").Value = 76924160
$ws.Range("I97").Value = 100000810
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 100000810
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -100000314
$ws.Range("N97").Value = -2992
$ws.Range("H110").Value = 1558.65
$ws.Range("I110").Value = 1509.1052
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 1509.1052
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = 535.8948
$ws.Range("N110").Value = -6590
$ws.Range("H116").Value = 942.3611
$ws.Range("I116").Value = 887.0357
$ws.Range("J116").Value = 1136
$ws.Range("K116").Value = 887.0357
$ws.Range("L116").Value = 1136
$ws.Range("M116").Value = 1406.9643
$ws.Range("N116").Value = -5724

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 942.3611
$ws.Range("I3").Value = 887.0357
$ws.Range("J3").Value = 1136
$ws.Range("K3").Value = 887.0357
$ws.Range("L3").Value = 1136
$ws.Range("M3").Value = -773.0357
$ws.Range("N3").Value = -1364
$ws.Range("H94").Value = 480.73334
$ws.Range("I94").Value = 400.84616
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 400.84616
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = 50.15384
$ws.Range("N94").Value = -1902
$ws.Range("H112").Value = 43746
$ws.Range("J112").Value = 43746
$ws.Range("L112").Value = 43746
$ws.Range("N112").Value = -46700

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1997.1666
$ws.Range("I22").Value = 329.55554
$ws.Range("J22").Value = 7000
$ws.Range("K22").Value = 329.55554
$ws.Range("L22").Value = 7000
$ws.Range("M22").Value = 20.44445999999999
$ws.Range("N22").Value = -7700
$ws.Range("H134").Value = 585139.4399999999
$ws.Range("I134").Value = 1308.1333
$ws.Range("J134").Value = 1558191.6
$ws.Range("K134").Value = 3924.3999
$ws.Range("L134").Value = 4674574.800000001
$ws.Range("M134").Value = -1389.3999
$ws.Range("N134").Value = -4679644.800000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 50815.5
$ws.Range("I130").Value = 60632
$ws.Range("J130").Value = 1733
$ws.Range("K130").Value = 181896
$ws.Range("L130").Value = 5199
$ws.Range("M130").Value = -176876
$ws.Range("N130").Value = -15239
$ws.Range("H134").Value = 23860934
$ws.Range("I134").Value = 31314852
$ws.Range("J134").Value = 8397.799999999999
$ws.Range("K134").Value = 93944556
$ws.Range("L134").Value = 25193.4
$ws.Range("M134").Value = -93939486
$ws.Range("N134").Value = -35333.39999999999
$ws.Range("H137").Value = 55567270
$ws.Range("J137").Value = 90925896
$ws.Range("L137").Value = 272777688
$ws.Range("N137").Value = -272787888

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 26000
$ws.Range("J26").Value = 26000
$ws.Range("L26").Value = 26000
$ws.Range("N26").Value = -26560
$ws.Range("H50").Value = 26000
$ws.Range("J50").Value = 26000
$ws.Range("L50").Value = 26000
$ws.Range("N50").Value = -26996

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3403.3076
$ws.Range("I7").Value = 2823.8
$ws.Range("J7").Value = 5335
$ws.Range("K7").Value = 2823.8
$ws.Range("L7").Value = 5335
$ws.Range("M7").Value = -2711.8
$ws.Range("N7").Value = -5559
$ws.Range("H22").Value = 872.3077
$ws.Range("I22").Value = 782.2222
$ws.Range("J22").Value = 1075
$ws.Range("K22").Value = 782.2222
$ws.Range("L22").Value = 1075
$ws.Range("M22").Value = -487.2222
$ws.Range("N22").Value = -1665
$ws.Range("H27").Value = 872.3077
$ws.Range("I27").Value = 782.2222
$ws.Range("J27").Value = 1075
$ws.Range("K27").Value = 782.2222
$ws.Range("L27").Value = 1075
$ws.Range("M27").Value = -675.2222
$ws.Range("N27").Value = -1289
$ws.Range("H40").Value = 2299.389
$ws.Range("I40").Value = 2183.3333
$ws.Range("J40").Value = 2531.5
$ws.Range("K40").Value = 2183.3333
$ws.Range("L40").Value = 2531.5
$ws.Range("M40").Value = -2047.3333
$ws.Range("N40").Value = -2803.5
$ws.Range("H46").Value = 3784.4285
$ws.Range("I46").Value = 1326.6666
$ws.Range("J46").Value = 4454.727
$ws.Range("K46").Value = 1326.6666
$ws.Range("L46").Value = 4454.727
$ws.Range("M46").Value = -1138.6666
$ws.Range("N46").Value = -4830.727
$ws.Range("H61").Value = 3981.6667
$ws.Range("I61").Value = 3780
$ws.Range("J61").Value = 4990
$ws.Range("K61").Value = 3780
$ws.Range("L61").Value = 4990
$ws.Range("M61").Value = -3578
$ws.Range("N61").Value = -5394
$ws.Range("H113").Value = 3981.6667
$ws.Range("I113").Value = 3780
$ws.Range("J113").Value = 4990
$ws.Range("K113").Value = 3780
$ws.Range("L113").Value = 4990
$ws.Range("M113").Value = -1610
$ws.Range("N113").Value = -9330
$ws.Range("H126").Value = 3403.3076
$ws.Range("I126").Value = 2823.8
$ws.Range("J126").Value = 5335
$ws.Range("K126").Value = 8471.400000000001
$ws.Range("L126").Value = 16005
$ws.Range("M126").Value = -6001.400000000001
$ws.Range("N126").Value = -20945
$ws.Range("H134").Value = 48749.75
$ws.Range("J134").Value = 48749.75
$ws.Range("L134").Value = 48749.75
$ws.Range("N134").Value = -58889.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 50000
$ws.Range("J26").Value = 50000
$ws.Range("L26").Value = 50000
$ws.Range("N26").Value = -50586
$ws.Range("H81").Value = 2044.5454
$ws.Range("I81").Value = 1598
$ws.Range("J81").Value = 2416.6667
$ws.Range("K81").Value = 3196
$ws.Range("L81").Value = 4833.3334
$ws.Range("M81").Value = -2135
$ws.Range("N81").Value = -6955.3334
$ws.Range("H84").Value = 2044.5454
$ws.Range("I84").Value = 1598
$ws.Range("J84").Value = 2416.6667
$ws.Range("K84").Value = 15980
$ws.Range("L84").Value = 24166.667
$ws.Range("M84").Value = -10676
$ws.Range("N84").Value = -34774.667
$ws.Range("H86").Value = 29999.5
$ws.Range("J86").Value = 29999.5
$ws.Range("L86").Value = 29999.5
$ws.Range("N86").Value = -32245.5
$ws.Range("H89").Value = 29999.5
$ws.Range("J89").Value = 29999.5
$ws.Range("L89").Value = 149997.5
$ws.Range("N89").Value = -161229.5
$ws.Range("H126").Value = 9804588
$ws.Range("I126").Value = 14706382
$ws.Range("J126").Value = 999
$ws.Range("K126").Value = 44119146
$ws.Range("L126").Value = 2997
$ws.Range("M126").Value = -44116676
$ws.Range("N126").Value = -7937
$ws.Range("H132").Value = 180842.75
$ws.Range("I132").Value = 251899.95
$ws.Range("K132").Value = 755699.8500000001
$ws.Range("M132").Value = -753169.8500000001
$ws.Range("H133").Value = 81855.8
$ws.Range("J133").Value = 81855.8
$ws.Range("L133").Value = 81855.8
$ws.Range("N133").Value = -91975.8
